$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 13.295
$ws.Range("C8").Value = -12.342
$ws.Range("C10").Value = -12.179
$ws.Range("C12").Value = -12.157
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("C18").Value = -12.157
$ws.Range("E20").Value = 12.761
$ws.Range("C25").Value = -12.324
